$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns A, B, E, F, G, I, M between row 2 and row 3
$cols = @("A", "B", "E", "F", "G", "I", "M")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $tmp = $cell2.Value2
    $cell2.Value = $cell3.Value2
    $cell3.Value = $tmp
}
